# display_template.xlsx edit:
#  - bump several "Shelf Length, m" values in column C (re-measured figures)
#  - drop the now-unused column G marker cells (G4, G15)
#  - extend the sheet's AutoFilter / _FilterDatabase range to cover all 104 rows
#  - register an extra (visible) _FilterDatabase_... defined name, as LibreOffice
#    does every time the autofilter range is redefined
#  - reset the view back to the top of the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Updated "Shelf Length, m" figures (column C) ------------------------
$ws.Range("C2").Value = 0.375
$ws.Range("C3").Value = 5.32
$ws.Range("C4").Value = 5.32
$ws.Range("C7").Value = 1.875
$ws.Range("C8").Value = 1.33
$ws.Range("C9").Value = 2.5
$ws.Range("C12").Value = 2.5
$ws.Range("C15").Value = 3.75
$ws.Range("C18").Value = 5.32
$ws.Range("C24").Value = 1.125
$ws.Range("C25").Value = 1.875
$ws.Range("C27").Value = 3.75
$ws.Range("C28").Value = 3.75
$ws.Range("C29").Value = 2.5
$ws.Range("C30").Value = 1.875
$ws.Range("C31").Value = 1.25
$ws.Range("C32").Value = 6
$ws.Range("C33").Value = 6
$ws.Range("C36").Value = 1.875
$ws.Range("C37").Value = 2.4
$ws.Range("C39").Value = 1.8
$ws.Range("C42").Value = 0.375
$ws.Range("C43").Value = 1.8
$ws.Range("C44").Value = 3.75
$ws.Range("C46").Value = 2.5
$ws.Range("C50").Value = 1.125
$ws.Range("C51").Value = 0.375
$ws.Range("C56").Value = 1.125
$ws.Range("C57").Value = 1.125
$ws.Range("C61").Value = 1.875
$ws.Range("C62").Value = 1.875
$ws.Range("C63").Value = 1.875
$ws.Range("C64").Value = 3.75
$ws.Range("C65").Value = 3.75
$ws.Range("C66").Value = 3.75
$ws.Range("C67").Value = 2.5
$ws.Range("C68").Value = 2.5
$ws.Range("C69").Value = 2.5
$ws.Range("C70").Value = 1.875
$ws.Range("C71").Value = 1.33
$ws.Range("C72").Value = 2.5
$ws.Range("C73").Value = 3.75
$ws.Range("C75").Value = 1.8
$ws.Range("C76").Value = 2.4
$ws.Range("C78").Value = 1.8
$ws.Range("C79").Value = 3.75
$ws.Range("C80").Value = 1.875
$ws.Range("C81").Value = 1.125
$ws.Range("C82").Value = 0.375
$ws.Range("C83").Value = 2.5
$ws.Range("C84").Value = 2.5
$ws.Range("C85").Value = 1.25
$ws.Range("C86").Value = 6
$ws.Range("C87").Value = 5.32
$ws.Range("C88").Value = 5.32
$ws.Range("C89").Value = 6
$ws.Range("C102").Value = 1.875
$ws.Range("C103").Value = 3.75
$ws.Range("C104").Value = 2.5

# ---- Drop the stray column-G cells (no longer used) -----------------------
$ws.Range("G4").Clear()
$ws.Range("G15").Clear()

# ---- Row 1 now wraps onto more lines -> taller header row -----------------
$ws.Rows.Item(1).RowHeight = 68.65

# ---- Re-apply the AutoFilter over the full (now 104-row) data range -------
$ws.AutoFilterMode = $false
$ws.Range("A1:F104").AutoFilter()

# ---- Keep the hidden _FilterDatabase name (localSheetId 0) in sync --------
foreach ($n in $wb.Names) {
  if ($n.Name -eq "Sheet1!_FilterDatabase" -and $n.Visible -eq $false) {
    $n.RefersTo = "=Sheet1!`$A`$1:`$F`$104"
  }
}
# LibreOffice stamps a fresh visible _FilterDatabase_... name every time the
# autofilter is (re)applied -- add the next one in the existing series.
$ws.Names.Add("_xlnm._FilterDatabase_0_0_0_0_0_0_0_0", "=Sheet1!`$A`$1:`$F`$22")

# ---- Reset the view back to the top-left of the sheet ----------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C19").Select()
